$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the description for RBLK (row 26) in column B
$ws.Range("B26").Value = "Tapon recibido"

# Update the selected cell to match the saved view state
$ws.Range("B27").Select()
